$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header row: add two new header cells P1, Q1 ---
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15

# Match the header formatting (bold font, thin border, centered/top align)
# used by the rest of row 1 (style carried by O1) onto the two new cells.
$ws.Range("O1").Copy()
$ws.Range("P1:Q1").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# --- Data rows 2-25 ---
for ($r = 2; $r -le 25; $r++) {
    # Flip the existing alternating 1/2 pattern in columns I, K, M, O
    $ws.Cells.Item($r, 9).Value  = 2   # I -> 2
    $ws.Cells.Item($r, 11).Value = 1   # K -> 1
    $ws.Cells.Item($r, 13).Value = 2   # M -> 2
    $ws.Cells.Item($r, 15).Value = 1   # O -> 1

    # New columns P, Q
    $ws.Cells.Item($r, 16).Value = 2   # P
    $ws.Cells.Item($r, 17).Value = 2   # Q
}
